$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 208, pushing the existing rows
# 208-214 down to 210-216.
$ws.Range("A208:A209").EntireRow.Insert()

# New row 208: Moscatel rosada
$ws.Cells.Item(208, 1).Value = 2
$ws.Cells.Item(208, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(208, 3).Value = "Coquimbo"
$ws.Cells.Item(208, 4).Value = 45021
$ws.Cells.Item(208, 5).Value = 4
$ws.Cells.Item(208, 6).Value = "Fruta"
$ws.Cells.Item(208, 7).Value = 100109
$ws.Cells.Item(208, 8).Value = "Uva"
$ws.Cells.Item(208, 9).Value = 100109001
$ws.Cells.Item(208, 10).Value = "Uva"
$ws.Cells.Item(208, 11).Value = "Moscatel rosada"
$ws.Cells.Item(208, 12).Value = "Primera"
$ws.Cells.Item(208, 13).Value = 700
$ws.Cells.Item(208, 14).Value = 12000
$ws.Cells.Item(208, 15).Value = 13000
$ws.Cells.Item(208, 16).Value = 12500
$ws.Cells.Item(208, 17).Value = "$/bandeja 12 kilos"
$ws.Cells.Item(208, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(208, 19).Value = 1042
$ws.Cells.Item(208, 20).Value = 12

# New row 209: Red Globe
$ws.Cells.Item(209, 1).Value = 2
$ws.Cells.Item(209, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 45021
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100109
$ws.Cells.Item(209, 8).Value = "Uva"
$ws.Cells.Item(209, 9).Value = 100109001
$ws.Cells.Item(209, 10).Value = "Uva"
$ws.Cells.Item(209, 11).Value = "Red Globe"
$ws.Cells.Item(209, 12).Value = "Primera"
$ws.Cells.Item(209, 13).Value = 400
$ws.Cells.Item(209, 14).Value = 7000
$ws.Cells.Item(209, 15).Value = 8000
$ws.Cells.Item(209, 16).Value = 7500
$ws.Cells.Item(209, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(209, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(209, 19).Value = 417
$ws.Cells.Item(209, 20).Value = 18
